$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.868.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.664.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.36%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.143.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.733.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.663.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "355.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("E24").Value = "  +10.00%  "
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "577.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.68%  "
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.36%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("E41").Value = "  +7.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("E50").Value = "  -5.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.821"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.80%  "
